$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row ranges (inclusive) in column C that receive a new constant Fitness value,
# per the commit "correction in sa algorithm and 746 logs".
$blocks = @(
    @{ Start = 2;  End = 7;   Value = 8244 },
    @{ Start = 8;  End = 12;  Value = 8081 },
    @{ Start = 13; End = 33;  Value = 7598 },
    @{ Start = 34; End = 252; Value = 7573 }
)

foreach ($block in $blocks) {
    $rangeAddr = "C" + $block.Start + ":C" + $block.End
    $ws.Range($rangeAddr).Value = $block.Value
}
